$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case with 380 kV done: updated bus voltage magnitude results (vm_pu)
# for rows 2-25 (bus indices 0-23), columns B-F and I-N.

$row = 2
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.00753248725116
$ws.Cells.Item($row, 4).Value = 1.032573116818584
$ws.Cells.Item($row, 5).Value = 1.010223396258179
$ws.Cells.Item($row, 6).Value = 1.019411235485953
$ws.Cells.Item($row, 9).Value = 1.0318397195802
$ws.Cells.Item($row, 10).Value = 1.012803578957126
$ws.Cells.Item($row, 11).Value = 1.035377958160746
$ws.Cells.Item($row, 12).Value = 1.013094236733674
$ws.Cells.Item($row, 13).Value = 1.022254562761938
$ws.Cells.Item($row, 14).Value = 1.008451638386829

$row = 3
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.008616901131992
$ws.Cells.Item($row, 4).Value = 1.03315344718361
$ws.Cells.Item($row, 5).Value = 1.011144145563421
$ws.Cells.Item($row, 6).Value = 1.021223612924085
$ws.Cells.Item($row, 9).Value = 1.032066743719867
$ws.Cells.Item($row, 10).Value = 1.013518791284481
$ws.Cells.Item($row, 11).Value = 1.035767700122445
$ws.Cells.Item($row, 12).Value = 1.013818686456051
$ws.Cells.Item($row, 13).Value = 1.023870194768746
$ws.Cells.Item($row, 14).Value = 1.008687483636365

$row = 4
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.009315478686534
$ws.Cells.Item($row, 4).Value = 1.033523199915388
$ws.Cells.Item($row, 5).Value = 1.011737666924019
$ws.Cells.Item($row, 6).Value = 1.022375100212514
$ws.Cells.Item($row, 9).Value = 1.032206753471097
$ws.Cells.Item($row, 10).Value = 1.013977996876173
$ws.Cells.Item($row, 11).Value = 1.036013444365753
$ws.Cells.Item($row, 12).Value = 1.014284577800729
$ws.Cells.Item($row, 13).Value = 1.024894016334502
$ws.Cells.Item($row, 14).Value = 1.008838887407205

$row = 5
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.009608423454855
$ws.Cells.Item($row, 4).Value = 1.03367726864701
$ws.Cells.Item($row, 5).Value = 1.011986646706227
$ws.Cells.Item($row, 6).Value = 1.022854136893522
$ws.Cells.Item($row, 9).Value = 1.032263967600344
$ws.Cells.Item($row, 10).Value = 1.014170194919927
$ws.Cells.Item($row, 11).Value = 1.036115215873595
$ws.Cells.Item($row, 12).Value = 1.01447975546556
$ws.Cells.Item($row, 13).Value = 1.025319293246648
$ws.Cells.Item($row, 14).Value = 1.008902251376798

$row = 6
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.00965756714841
$ws.Cells.Item($row, 4).Value = 1.033703056910239
$ws.Cells.Item($row, 5).Value = 1.012028420198448
$ws.Cells.Item($row, 6).Value = 1.022934274236596
$ws.Cells.Item($row, 9).Value = 1.032273477676487
$ws.Cells.Item($row, 10).Value = 1.014202416028223
$ws.Cells.Item($row, 11).Value = 1.036132213581076
$ws.Cells.Item($row, 12).Value = 1.01451248675224
$ws.Cells.Item($row, 13).Value = 1.025390398767039
$ws.Cells.Item($row, 14).Value = 1.00891287374132

$row = 7
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.009319395915563
$ws.Cells.Item($row, 4).Value = 1.033525263988517
$ws.Cells.Item($row, 5).Value = 1.011740995904858
$ws.Cells.Item($row, 6).Value = 1.022381520913021
$ws.Cells.Item($row, 9).Value = 1.032207524431952
$ws.Cells.Item($row, 10).Value = 1.013980568375273
$ws.Cells.Item($row, 11).Value = 1.036014810285556
$ws.Cells.Item($row, 12).Value = 1.014287188450798
$ws.Cells.Item($row, 13).Value = 1.024899719048859
$ws.Cells.Item($row, 14).Value = 1.008839735201407

$row = 8
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.007899618740386
$ws.Cells.Item($row, 4).Value = 1.032770437129139
$ws.Cells.Item($row, 5).Value = 1.01053504030171
$ws.Cells.Item($row, 6).Value = 1.020028158164776
$ws.Cells.Item($row, 9).Value = 1.03191787236856
$ws.Cells.Item($row, 10).Value = 1.013046034379705
$ws.Cells.Item($row, 11).Value = 1.035511010037386
$ws.Cells.Item($row, 12).Value = 1.013339666834163
$ws.Cells.Item($row, 13).Value = 1.022805066567785
$ws.Cells.Item($row, 14).Value = 1.008531593931909

$row = 9
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.005373606920761
$ws.Cells.Item($row, 4).Value = 1.031396031062528
$ws.Cells.Item($row, 5).Value = 1.008392371870897
$ws.Cells.Item($row, 6).Value = 1.015716767510367
$ws.Cells.Item($row, 9).Value = 1.031354511930706
$ws.Cells.Item($row, 10).Value = 1.011371512972327
$ws.Cells.Item($row, 11).Value = 1.034573701087805
$ws.Cells.Item($row, 12).Value = 1.011647720046818
$ws.Cells.Item($row, 13).Value = 1.01894696320431
$ws.Cells.Item($row, 14).Value = 1.007979290875204

$row = 10
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.003672803897405
$ws.Cells.Item($row, 4).Value = 1.030449659063463
$ws.Cells.Item($row, 5).Value = 1.006951666841062
$ws.Cells.Item($row, 6).Value = 1.012729267473493
$ws.Cells.Item($row, 9).Value = 1.030943028294413
$ws.Cells.Item($row, 10).Value = 1.010236053948717
$ws.Cells.Item($row, 11).Value = 1.033915206272555
$ws.Cells.Item($row, 12).Value = 1.010504366685369
$ws.Cells.Item($row, 13).Value = 1.016260123399946
$ws.Cells.Item($row, 14).Value = 1.007604674811571

$row = 11
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002932224655331
$ws.Cells.Item($row, 4).Value = 1.030032647105261
$ws.Cells.Item($row, 5).Value = 1.006324817465804
$ws.Cells.Item($row, 6).Value = 1.011408144054811
$ws.Cells.Item($row, 9).Value = 1.03075625361507
$ws.Cells.Item($row, 10).Value = 1.009739744587624
$ws.Cells.Item($row, 11).Value = 1.033622014523531
$ws.Cells.Item($row, 12).Value = 1.010005534450627
$ws.Cells.Item($row, 13).Value = 1.015068873702881
$ws.Cells.Item($row, 14).Value = 1.007440904412309

$row = 12
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002656509699459
$ws.Cells.Item($row, 4).Value = 1.029876656505176
$ws.Cells.Item($row, 5).Value = 1.006091515535646
$ws.Cells.Item($row, 6).Value = 1.010913227018937
$ws.Cells.Item($row, 9).Value = 1.03068557737884
$ws.Cells.Item($row, 10).Value = 1.009554684715342
$ws.Cells.Item($row, 11).Value = 1.033511891315245
$ws.Cells.Item($row, 12).Value = 1.009819673153888
$ws.Cells.Item($row, 13).Value = 1.014622154482272
$ws.Cells.Item($row, 14).Value = 1.007379835210691

$row = 13
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002715680237042
$ws.Cells.Item($row, 4).Value = 1.029910166651713
$ws.Cells.Item($row, 5).Value = 1.006141580659501
$ws.Cells.Item($row, 6).Value = 1.011019579125553
$ws.Cells.Item($row, 9).Value = 1.030700796622887
$ws.Cells.Item($row, 10).Value = 1.009594412915329
$ws.Cells.Item($row, 11).Value = 1.033535568402463
$ws.Cells.Item($row, 12).Value = 1.009859567102926
$ws.Cells.Item($row, 13).Value = 1.014718169875498
$ws.Cells.Item($row, 14).Value = 1.00739294556844

$row = 14
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002909446910748
$ws.Cells.Item($row, 4).Value = 1.030019775235722
$ws.Cells.Item($row, 5).Value = 1.006305542138783
$ws.Cells.Item($row, 6).Value = 1.01136731994632
$ws.Cells.Item($row, 9).Value = 1.030750438058509
$ws.Cells.Item($row, 10).Value = 1.009724461994317
$ws.Cells.Item($row, 11).Value = 1.033612936615005
$ws.Cells.Item($row, 12).Value = 1.009990182836506
$ws.Cells.Item($row, 13).Value = 1.015032034495197
$ws.Cells.Item($row, 14).Value = 1.007435861279637

$row = 15
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.003028748996229
$ws.Cells.Item($row, 4).Value = 1.030087163511143
$ws.Cells.Item($row, 5).Value = 1.006406502691489
$ws.Cells.Item($row, 6).Value = 1.011581017115706
$ws.Cells.Item($row, 9).Value = 1.030780851303915
$ws.Cells.Item($row, 10).Value = 1.00980449533094
$ws.Cells.Item($row, 11).Value = 1.033660443980835
$ws.Cells.Item($row, 12).Value = 1.010070583322728
$ws.Cells.Item($row, 13).Value = 1.01522485383715
$ws.Cells.Item($row, 14).Value = 1.007462271480774

$row = 16
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.003721865896529
$ws.Cells.Item($row, 4).Value = 1.030477181843028
$ws.Cells.Item($row, 5).Value = 1.006993204476967
$ws.Cells.Item($row, 6).Value = 1.012816361385507
$ws.Cells.Item($row, 9).Value = 1.030955242079252
$ws.Cells.Item($row, 10).Value = 1.010268893582022
$ws.Cells.Item($row, 11).Value = 1.033934493983286
$ws.Cells.Item($row, 12).Value = 1.010537392716253
$ws.Cells.Item($row, 13).Value = 1.016338591752623
$ws.Cells.Item($row, 14).Value = 1.007615510582986

$row = 17
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.004155528526728
$ws.Cells.Item($row, 4).Value = 1.030719889773349
$ws.Cells.Item($row, 5).Value = 1.007360413535873
$ws.Cells.Item($row, 6).Value = 1.013583852677949
$ws.Cells.Item($row, 9).Value = 1.031062325142212
$ws.Cells.Item($row, 10).Value = 1.01055894690869
$ws.Cells.Item($row, 11).Value = 1.034104235244844
$ws.Cells.Item($row, 12).Value = 1.010829199015409
$ws.Cells.Item($row, 13).Value = 1.017029721374806
$ws.Cells.Item($row, 14).Value = 1.007711213689684

$row = 18
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.004408079883137
$ws.Cells.Item($row, 4).Value = 1.030860760521346
$ws.Cells.Item($row, 5).Value = 1.00757431015627
$ws.Cells.Item($row, 6).Value = 1.014028866608065
$ws.Cells.Item($row, 9).Value = 1.031123955714785
$ws.Cells.Item($row, 10).Value = 1.010727682348349
$ws.Cells.Item($row, 11).Value = 1.034202465434769
$ws.Cells.Item($row, 12).Value = 1.010999043248283
$ws.Cells.Item($row, 13).Value = 1.017430164241435
$ws.Cells.Item($row, 14).Value = 1.007766885490669

$row = 19
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.004494126426985
$ws.Cells.Item($row, 4).Value = 1.030908675869869
$ws.Cells.Item($row, 5).Value = 1.00764719442329
$ws.Cells.Item($row, 6).Value = 1.014180156978782
$ws.Cells.Item($row, 9).Value = 1.031144829732427
$ws.Cells.Item($row, 10).Value = 1.010785141123359
$ws.Cells.Item($row, 11).Value = 1.034235827825111
$ws.Cells.Item($row, 12).Value = 1.011056894699554
$ws.Cells.Item($row, 13).Value = 1.017566251718437
$ws.Cells.Item($row, 14).Value = 1.00778584275726

$row = 20
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.004109041749097
$ws.Cells.Item($row, 4).Value = 1.03069392162693
$ws.Cells.Item($row, 5).Value = 1.007321045562579
$ws.Cells.Item($row, 6).Value = 1.013501782723244
$ws.Cells.Item($row, 9).Value = 1.031050921962255
$ws.Cells.Item($row, 10).Value = 1.010527873317441
$ws.Cells.Item($row, 11).Value = 1.034086104039902
$ws.Cells.Item($row, 12).Value = 1.010797928410777
$ws.Cells.Item($row, 13).Value = 1.016955847439433
$ws.Cells.Item($row, 14).Value = 1.007700961211788

$row = 21
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002852404967144
$ws.Cells.Item($row, 4).Value = 1.029987528501712
$ws.Cells.Item($row, 5).Value = 1.006257272404384
$ws.Cells.Item($row, 6).Value = 1.011265035177111
$ws.Cells.Item($row, 9).Value = 1.030735855839868
$ws.Cells.Item($row, 10).Value = 1.009686185423552
$ws.Cells.Item($row, 11).Value = 1.033590187314885
$ws.Cells.Item($row, 12).Value = 1.009951735644976
$ws.Cells.Item($row, 13).Value = 1.014939726578159
$ws.Cells.Item($row, 14).Value = 1.007423230258451

$row = 22
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002058651370121
$ws.Cells.Item($row, 4).Value = 1.029537058293506
$ws.Cells.Item($row, 5).Value = 1.005585757138569
$ws.Cells.Item($row, 6).Value = 1.009834412521431
$ws.Cells.Item($row, 9).Value = 1.030530236667116
$ws.Cells.Item($row, 10).Value = 1.009152878387606
$ws.Cells.Item($row, 11).Value = 1.033271329140068
$ws.Cells.Item($row, 12).Value = 1.009416381856719
$ws.Cells.Item($row, 13).Value = 1.013647573607611
$ws.Cells.Item($row, 14).Value = 1.007247233420732

$row = 23
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.002479785828946
$ws.Cells.Item($row, 4).Value = 1.029776464218435
$ws.Cells.Item($row, 5).Value = 1.005941997227995
$ws.Cells.Item($row, 6).Value = 1.010595135633201
$ws.Cells.Item($row, 9).Value = 1.030639955301891
$ws.Cells.Item($row, 10).Value = 1.009435987122854
$ws.Cells.Item($row, 11).Value = 1.033441033506493
$ws.Cells.Item($row, 12).Value = 1.009700500742005
$ws.Cells.Item($row, 13).Value = 1.014334913337438
$ws.Cells.Item($row, 14).Value = 1.007340664297542

$row = 24
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.004130048348609
$ws.Cells.Item($row, 4).Value = 1.030705657667191
$ws.Cells.Item($row, 5).Value = 1.007338835151798
$ws.Cells.Item($row, 6).Value = 1.013538874791295
$ws.Cells.Item($row, 9).Value = 1.031056077130349
$ws.Cells.Item($row, 10).Value = 1.010541915519125
$ws.Cells.Item($row, 11).Value = 1.034094299151806
$ws.Cells.Item($row, 12).Value = 1.010812059367685
$ws.Cells.Item($row, 13).Value = 1.01698923617389
$ws.Cells.Item($row, 14).Value = 1.007705594329287

$row = 25
$ws.Cells.Item($row, 2).Value = 1.02
$ws.Cells.Item($row, 3).Value = 1.006029560573973
$ws.Cells.Item($row, 4).Value = 1.031756628232588
$ws.Cells.Item($row, 5).Value = 1.008948431945024
$ws.Cells.Item($row, 6).Value = 1.016851118954065
$ws.Cells.Item($row, 9).Value = 1.031506458879438
$ws.Cells.Item($row, 10).Value = 1.011807750445489
$ws.Cells.Item($row, 11).Value = 1.034821919139515
$ws.Cells.Item($row, 12).Value = 1.012087810248591
$ws.Cells.Item($row, 13).Value = 1.019964414708994
$ws.Cells.Item($row, 14).Value = 1.008123193473792

